# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets,
# matching the refreshed data output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 20
$ws1.Range("F4").Value = 1374
$ws1.Range("F5").Value = 312
$ws1.Range("F7").Value = 10669
$ws1.Range("F12").Value = 701
$ws1.Range("F13").Value = 12040
$ws1.Range("F14").Value = 12482
$ws1.Range("F16").Value = 120
$ws1.Range("F19").Value = 75

# Sheet "全部类型" (all types, combined feed)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 20
$ws4.Range("F5").Value = 1374
$ws4.Range("F6").Value = 312
$ws4.Range("F8").Value = 10669
$ws4.Range("F13").Value = 701
$ws4.Range("F14").Value = 12040
$ws4.Range("F15").Value = 12482
$ws4.Range("F17").Value = 120
$ws4.Range("F20").Value = 75
